$d = $word.ActiveDocument

# Locate the text "JavaScript med DOM " (without the ampersand) so we can
# find the character offset right before the "&" that follows it.
$find = $d.Content
$find.Find.ClearFormatting()
$find.Find.Execute("JavaScript med DOM ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Re-adding the "_GoBack" bookmark at the new location moves it there
# (splitting the run "JavaScript med DOM &" into "JavaScript med DOM "
# and "&" in the process), matching the target edit.
$d.Bookmarks.Add("_GoBack", $d.Range($find.End, $find.End))
